$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3172.1765
$ws.Range("I113").Value = 2199.8
$ws.Range("J113").Value = 3577.3333
$ws.Range("K113").Value = 2199.8
$ws.Range("L113").Value = 3577.3333
$ws.Range("M113").Value = 1054.2
$ws.Range("N113").Value = -10085.3333
$ws.Range("H123").Value = 54990
$ws.Range("J123").Value = 54990
$ws.Range("L123").Value = 54990
$ws.Range("N123").Value = -64790
$ws.Range("H124").Value = 78800
$ws.Range("J124").Value = 78800
$ws.Range("L124").Value = 78800
$ws.Range("N124").Value = -88620
$ws.Range("H125").Value = 9854.666999999999
$ws.Range("I125").Value = 2074.8333
$ws.Range("J125").Value = 17634.5
$ws.Range("K125").Value = 18673.4997
$ws.Range("L125").Value = 158710.5
$ws.Range("M125").Value = -16213.4997
$ws.Range("N125").Value = -163630.5
$ws.Range("H129").Value = 1042.1428
$ws.Range("I129").Value = 900
$ws.Range("J129").Value = 1046.3235
$ws.Range("K129").Value = 2700
$ws.Range("L129").Value = 3138.9705
$ws.Range("M129").Value = 2300
$ws.Range("N129").Value = -13138.9705
$ws.Range("H132").Value = 3239.8
$ws.Range("I132").Value = 2683.9473
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8051.841899999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5521.841899999999
$ws.Range("N132").Value = -20060
$ws.Range("H137").Value = 536486.0600000001
$ws.Range("I137").Value = 1447.8485
$ws.Range("J137").Value = 1465762.9
$ws.Range("K137").Value = 4343.5455
$ws.Range("L137").Value = 4397288.699999999
$ws.Range("M137").Value = -1793.5455
$ws.Range("N137").Value = -4402388.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1212.3438
$ws.Range("I2").Value = 1191.8
$ws.Range("J2").Value = 1285.7142
$ws.Range("K2").Value = 1191.8
$ws.Range("L2").Value = 1285.7142
$ws.Range("M2").Value = -1078.8
$ws.Range("N2").Value = -1511.7142
$ws.Range("H45").Value = 1965.5385
$ws.Range("I45").Value = 1795
$ws.Range("J45").Value = 2903.5
$ws.Range("K45").Value = 1795
$ws.Range("L45").Value = 2903.5
$ws.Range("M45").Value = -1418
$ws.Range("N45").Value = -3657.5
$ws.Range("H61").Value = 6560.5356
$ws.Range("I61").Value = 5092.1
$ws.Range("J61").Value = 10231.625
$ws.Range("K61").Value = 5092.1
$ws.Range("L61").Value = 10231.625
$ws.Range("M61").Value = -4880.1
$ws.Range("N61").Value = -10655.625
$ws.Range("H116").Value = 1212.3438
$ws.Range("I116").Value = 1191.8
$ws.Range("J116").Value = 1285.7142
$ws.Range("K116").Value = 1191.8
$ws.Range("L116").Value = 1285.7142
$ws.Range("M116").Value = 1102.2
$ws.Range("N116").Value = -5873.7142
$ws.Range("H122").Value = 4034876.2
$ws.Range("I122").Value = 3457.4614
$ws.Range("J122").Value = 6946456.5
$ws.Range("K122").Value = 10372.3842
$ws.Range("L122").Value = 20839369.5
$ws.Range("M122").Value = -7922.3842
$ws.Range("N122").Value = -20844269.5
$ws.Range("H136").Value = 6560.5356
$ws.Range("I136").Value = 5092.1
$ws.Range("J136").Value = 10231.625
$ws.Range("K136").Value = 15276.3
$ws.Range("L136").Value = 30694.875
$ws.Range("M136").Value = -12726.3
$ws.Range("N136").Value = -35794.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1212.3438
$ws.Range("I3").Value = 1191.8
$ws.Range("J3").Value = 1285.7142
$ws.Range("K3").Value = 1191.8
$ws.Range("L3").Value = 1285.7142
$ws.Range("M3").Value = -1077.8
$ws.Range("N3").Value = -1513.7142
$ws.Range("H118").Value = 57179.875
$ws.Range("J118").Value = 57179.875
$ws.Range("L118").Value = 57179.875
$ws.Range("N118").Value = -60493.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2539.5652
$ws.Range("I31").Value = 1953.6571
$ws.Range("J31").Value = 4403.8184
$ws.Range("K31").Value = 1953.6571
$ws.Range("L31").Value = 4403.8184
$ws.Range("M31").Value = -1658.6571
$ws.Range("N31").Value = -4993.8184
$ws.Range("H34").Value = 2539.5652
$ws.Range("I34").Value = 1953.6571
$ws.Range("J34").Value = 4403.8184
$ws.Range("K34").Value = 1953.6571
$ws.Range("L34").Value = 4403.8184
$ws.Range("M34").Value = -1751.6571
$ws.Range("N34").Value = -4807.8184
$ws.Range("H58").Value = 2167289.5
$ws.Range("I58").Value = 3638130.5
$ws.Range("K58").Value = 3638130.5
$ws.Range("M58").Value = -3637927.5
$ws.Range("H86").Value = 1918
$ws.Range("I86").Value = 1200
$ws.Range("K86").Value = 1200
$ws.Range("M86").Value = -77
$ws.Range("H89").Value = 1918
$ws.Range("I89").Value = 1200
$ws.Range("K89").Value = 6000
$ws.Range("M89").Value = -384
$ws.Range("H94").Value = 1557.5333
$ws.Range("I94").Value = 1984.8
$ws.Range("J94").Value = 1343.9
$ws.Range("K94").Value = 1984.8
$ws.Range("L94").Value = 1343.9
$ws.Range("M94").Value = -1533.8
$ws.Range("N94").Value = -2245.9
$ws.Range("H107").Value = 752.8461
$ws.Range("I107").Value = 870.75
$ws.Range("K107").Value = 870.75
$ws.Range("M107").Value = 1049.25
$ws.Range("H122").Value = 6356.4
$ws.Range("I122").Value = 3320.375
$ws.Range("J122").Value = 11753.777
$ws.Range("K122").Value = 9961.125
$ws.Range("L122").Value = 35261.331
$ws.Range("M122").Value = -7511.125
$ws.Range("N122").Value = -40161.331
$ws.Range("H132").Value = 2026.6487
$ws.Range("I132").Value = 1542.5862
$ws.Range("J132").Value = 3781.375
$ws.Range("K132").Value = 4627.7586
$ws.Range("L132").Value = 11344.125
$ws.Range("M132").Value = -2097.7586
$ws.Range("N132").Value = -16404.125
$ws.Range("H134").Value = 2921.3052
$ws.Range("I134").Value = 1972.4722
$ws.Range("K134").Value = 5917.4166
$ws.Range("M134").Value = -3382.4166
$ws.Range("H136").Value = 2167289.5
$ws.Range("I136").Value = 3638130.5
$ws.Range("K136").Value = 10914391.5
$ws.Range("M136").Value = -10911841.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1275.75
$ws.Range("J92").Value = 1267.6666
$ws.Range("L92").Value = 3802.9998
$ws.Range("N92").Value = -6298.9998
$ws.Range("H107").Value = 1463.5294
$ws.Range("I107").Value = 365.7143
$ws.Range("J107").Value = 1748.1482
$ws.Range("K107").Value = 1097.1429
$ws.Range("L107").Value = 5244.444600000001
$ws.Range("M107").Value = 822.8571000000002
$ws.Range("N107").Value = -9084.444600000001
$ws.Range("H122").Value = 601.0417
$ws.Range("I122").Value = 459.85715
$ws.Range("J122").Value = 798.7
$ws.Range("K122").Value = 4138.71435
$ws.Range("L122").Value = 7188.3
$ws.Range("M122").Value = -1688.71435
$ws.Range("N122").Value = -12088.3
$ws.Range("H132").Value = 1621.375
$ws.Range("I132").Value = 1920.4
$ws.Range("J132").Value = 1407.7858
$ws.Range("K132").Value = 17283.6
$ws.Range("L132").Value = 12670.0722
$ws.Range("M132").Value = -14753.6
$ws.Range("N132").Value = -17730.0722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3042
$ws.Range("I102").Value = 3585
$ws.Range("J102").Value = 1413
$ws.Range("K102").Value = 3585
$ws.Range("L102").Value = 1413
$ws.Range("M102").Value = -1963
$ws.Range("N102").Value = -4657
$ws.Range("H113").Value = 2443.375
$ws.Range("I113").Value = 2046.8889
$ws.Range("J113").Value = 2953.1428
$ws.Range("K113").Value = 2046.8889
$ws.Range("L113").Value = 2953.1428
$ws.Range("M113").Value = 123.1111000000001
$ws.Range("N113").Value = -7293.1428
$ws.Range("H122").Value = 4641.3125
$ws.Range("I122").Value = 6114.636
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 18343.908
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -15893.908
$ws.Range("N122").Value = -9100
$ws.Range("H126").Value = 2733.4783
$ws.Range("I126").Value = 1489.1666
$ws.Range("J126").Value = 4090.9092
$ws.Range("K126").Value = 4467.4998
$ws.Range("L126").Value = 12272.7276
$ws.Range("M126").Value = -1997.4998
$ws.Range("N126").Value = -17212.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4296.593
$ws.Range("I40").Value = 4000.3809
$ws.Range("J40").Value = 5333.3335
$ws.Range("K40").Value = 4000.3809
$ws.Range("L40").Value = 5333.3335
$ws.Range("M40").Value = -3864.3809
$ws.Range("N40").Value = -5605.3335
$ws.Range("H82").Value = 1741.9231
$ws.Range("J82").Value = 2350.6
$ws.Range("L82").Value = 2350.6
$ws.Range("N82").Value = -3072.6
$ws.Range("H85").Value = 1741.9231
$ws.Range("J85").Value = 2350.6
$ws.Range("L85").Value = 2350.6
$ws.Range("N85").Value = -4846.6
$ws.Range("H93").Value = 1500
$ws.Range("I93").Value = 1500
$ws.Range("K93").Value = 1500
$ws.Range("M93").Value = -252

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1655.5555
$ws.Range("I81").Value = 983.3333
$ws.Range("K81").Value = 1966.6666
$ws.Range("M81").Value = -905.6666
$ws.Range("H84").Value = 1655.5555
$ws.Range("I84").Value = 983.3333
$ws.Range("K84").Value = 9833.333000000001
$ws.Range("M84").Value = -4529.333000000001
$ws.Range("H113").Value = 1166.5
$ws.Range("I113").Value = 558.0909
$ws.Range("J113").Value = 1518.7368
$ws.Range("K113").Value = 1674.2727
$ws.Range("L113").Value = 4556.2104
$ws.Range("M113").Value = 495.7273
$ws.Range("N113").Value = -8896.2104
$ws.Range("H132").Value = 1701.8214
$ws.Range("I132").Value = 1682.4762
$ws.Range("J132").Value = 1759.8572
$ws.Range("K132").Value = 5047.4286
$ws.Range("L132").Value = 5279.571599999999
$ws.Range("M132").Value = -2517.4286
$ws.Range("N132").Value = -10339.5716
$ws.Range("H136").Value = 5159.021
$ws.Range("I136").Value = 2789.3635
$ws.Range("J136").Value = 7164.115
$ws.Range("K136").Value = 8368.0905
$ws.Range("L136").Value = 21492.345
$ws.Range("M136").Value = -5818.0905
$ws.Range("N136").Value = -26592.345
